# edit.ps1 - apply the two grammar-rule fixes described by the diff:
#   1. "V -> i" paragraph: mark the bare "i" run as spell-checked
#      (wrap it in w:proofErr spellStart/spellEnd, after splitting the
#      trailing " i" run into " " + "i").
#   2. "S -> repeat S until B" paragraph: the rule was missing the
#      REPEAT_N nonterminal. Rewrite it as "S -> REPEAT_N S until B"
#      and add a new paragraph "REPEAT_N -> repeat" right after it,
#      matching the WHILE_D/WHILE_N pattern used elsewhere in the doc.

$d = $word.ActiveDocument

$pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$pkgFooter = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

function Set-ParagraphXml($paragraph, [string]$innerWpXml) {
  $paragraph.Range.InsertXML($pkgHeader + '<w:body>' + $innerWpXml + '</w:body>' + $pkgFooter)
}

# ---------------------------------------------------------------------
# Change 1: "V -> i"  =>  "V -> " + proofErr-wrapped "i"
# ---------------------------------------------------------------------
$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
  $cand = $d.Paragraphs.Item($i)
  if ($cand.Range.Text -eq "V -> i`r") {
    $found = $true
    $vp1 = '<w:p><w:pPr><w:ind w:firstLine="420"/></w:pPr>' +
      '<w:r><w:t xml:space="preserve">V </w:t></w:r>' +
      '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>-&gt;</w:t></w:r>' +
      '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
      '<w:proofErr w:type="spellStart"/>' +
      '<w:r><w:t>i</w:t></w:r>' +
      '<w:proofErr w:type="spellEnd"/>' +
      '</w:p>'
    Set-ParagraphXml $cand $vp1
    break
  }
}
if (-not $found) {
  Write-Output "WARNING: 'V -> i' paragraph not found"
}

# ---------------------------------------------------------------------
# Change 2: "S -> repeat S until B"
#   =>  "S -> REPEAT_N S until B"
#       + new paragraph "REPEAT_N -> repeat"
# ---------------------------------------------------------------------
$found2 = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
  $cand = $d.Paragraphs.Item($i)
  if ($cand.Range.Text -eq "S -> repeat S until B`r") {
    $found2 = $true

    $sp1 = '<w:p><w:pPr><w:ind w:firstLineChars="200" w:firstLine="420"/></w:pPr>' +
      '<w:r><w:t>S</w:t></w:r>' +
      '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
      '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>-&gt;</w:t></w:r>' +
      '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
      '<w:r><w:t>REPEAT_N</w:t></w:r>' +
      '<w:r><w:t xml:space="preserve"> S </w:t></w:r>' +
      '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>until</w:t></w:r>' +
      '<w:r><w:t xml:space="preserve"> B</w:t></w:r>' +
      '</w:p>'
    Set-ParagraphXml $cand $sp1

    # Re-fetch (index is stable; content was replaced in place) and
    # add the new "REPEAT_N -> repeat" paragraph right after it.
    $cand = $d.Paragraphs.Item($i)
    $cand.Range.InsertParagraphAfter()

    $newPara = $d.Paragraphs.Item($i + 1)
    $sp2 = '<w:p><w:pPr><w:ind w:firstLineChars="200" w:firstLine="420"/></w:pPr>' +
      '<w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/><w:lang w:eastAsia="zh-TW"/></w:rPr>' +
      '<w:t xml:space="preserve">REPEAT_N -&gt; </w:t></w:r>' +
      '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>repeat</w:t></w:r>' +
      '</w:p>'
    Set-ParagraphXml $newPara $sp2

    break
  }
}
if (-not $found2) {
  Write-Output "WARNING: 'S -> repeat S until B' paragraph not found"
}

Write-Output "Done. found1=$found found2=$found2 paragraphs=$($d.Paragraphs.Count)"
